$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: update status text for both locales
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet: status text, new Target/Handback hyperlinks, new handback datetime
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B2").Value = "Handed back: in sync with en-US"
$wsZh.Range("B3").Value = "Handed back: in sync with en-US"

$zhMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/2fdb5487e30d3d51a75517f3788e29b775d102ef/e2e/010f252c-d685-4920-96f4-db38c66a4380.md"
$zhXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2925f8078bf6df8a04fb499d02adbccabfade8cc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/010f252c-d685-4920-96f4-db38c66a4380.a6349a455b2828e633f96fd517ca981586124ea5.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $zhMdUrl, "", "", "010f252c-d685-4920-96f4-db38c66a4380.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhXlfUrl2, "", "", "010f252c-d685-4920-96f4-db38c66a4380.a6349a455b2828e633f96fd517ca981586124ea5.zh-cn.xlf")
$wsZh.Range("E2").Style = $wsZh.Range("A2").Style
$wsZh.Range("F2").Style = $wsZh.Range("A2").Style

$zhMdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/2fdb5487e30d3d51a75517f3788e29b775d102ef/e2e/ef223ad3-f544-4de1-8332-db0553c40196.md"
$zhXlfUrl3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2925f8078bf6df8a04fb499d02adbccabfade8cc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ef223ad3-f544-4de1-8332-db0553c40196.cba60d90b494e348eeb5cbd672abba22f553ffb5.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), $zhMdUrl2, "", "", "ef223ad3-f544-4de1-8332-db0553c40196.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhXlfUrl3, "", "", "ef223ad3-f544-4de1-8332-db0553c40196.cba60d90b494e348eeb5cbd672abba22f553ffb5.zh-cn.xlf")
$wsZh.Range("E3").Style = $wsZh.Range("A2").Style
$wsZh.Range("F3").Style = $wsZh.Range("A2").Style

$wsZh.Range("G2").Value = "2016-03-04 06:36:30"
$wsZh.Range("G3").Value = "2016-03-04 06:36:30"

# ---------------------------------------------------------------------
# de-de sheet: status text, new Target/Handback hyperlinks, new handback datetime
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B2").Value = "Handed back: in sync with en-US"
$wsDe.Range("B3").Value = "Handed back: in sync with en-US"

$deMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/2fdb5487e30d3d51a75517f3788e29b775d102ef/e2e/010f252c-d685-4920-96f4-db38c66a4380.md"
$deXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f86b44c8df6653577e7bfcdc76ba80260f1941d0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/010f252c-d685-4920-96f4-db38c66a4380.a6349a455b2828e633f96fd517ca981586124ea5.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $deMdUrl, "", "", "010f252c-d685-4920-96f4-db38c66a4380.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deXlfUrl2, "", "", "010f252c-d685-4920-96f4-db38c66a4380.a6349a455b2828e633f96fd517ca981586124ea5.de-de.xlf")
$wsDe.Range("E2").Style = $wsDe.Range("A2").Style
$wsDe.Range("F2").Style = $wsDe.Range("A2").Style

$deMdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/2fdb5487e30d3d51a75517f3788e29b775d102ef/e2e/ef223ad3-f544-4de1-8332-db0553c40196.md"
$deXlfUrl3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f86b44c8df6653577e7bfcdc76ba80260f1941d0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ef223ad3-f544-4de1-8332-db0553c40196.cba60d90b494e348eeb5cbd672abba22f553ffb5.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), $deMdUrl2, "", "", "ef223ad3-f544-4de1-8332-db0553c40196.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deXlfUrl3, "", "", "ef223ad3-f544-4de1-8332-db0553c40196.cba60d90b494e348eeb5cbd672abba22f553ffb5.de-de.xlf")
$wsDe.Range("E3").Style = $wsDe.Range("A2").Style
$wsDe.Range("F3").Style = $wsDe.Range("A2").Style

$wsDe.Range("G2").Value = "2016-03-04 06:36:57"
$wsDe.Range("G3").Value = "2016-03-04 06:36:57"
